# DOMA-3372: add the `role` field to the contacts excel export template.
#
# The template has a header row (i18n placeholders), a "row i" sample row
# and a "row i+1" sample row, each with columns for name / address /
# unitName / unitType / phone / email. We add a new trailing "role" column
# (G) with the matching i18n / contacts[i] / contacts[i+1] placeholders,
# copying the look (style/format/width) of the existing last column (F).
# The old unused "empty body" rows 4-10 (kept around in the source sheet
# only as leftover style placeholders) are removed as part of the cleanup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "role" column (G), matching column F's formatting ---
$ws.Range("F1:F3").Copy()
$ws.Range("G1:G3").PasteSpecial(-4122) # xlPasteFormats
$ws.Columns("G").ColumnWidth = $ws.Columns("F").ColumnWidth

$ws.Range("G1").Value = "{d.i18n.role}"
$ws.Range("G2").Value = "{d.contacts[i].role}"
$ws.Range("G3").Value = "{d.contacts[i+1].role}"

# --- Drop the leftover empty rows 4-10 ---
$ws.Range("A4:G10").EntireRow.Delete()
